$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.372.93"
$ws.Range("E2").Value = "  +0.35%  "

$ws.Range("D3").Value = "2.650.59"
$ws.Range("E3").Value = "  +0.61%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "520.84"
$ws.Range("E5").Value = "  +0.92%  "

$ws.Range("D6").Value = "146.77"
$ws.Range("E6").Value = "  +0.94%  "

$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").Value = "2.663.43"
$ws.Range("E9").Value = "  +0.98%  "

$ws.Range("E10").Value = "  +9.81%  "

$ws.Range("D11").Value = "0.103"
$ws.Range("E11").Value = "  -1.97%  "

$ws.Range("E12").Value = "  -0.15%  "

$ws.Range("E13").Value = "  +1.96%  "

$ws.Range("D14").Value = "3.114.72"
$ws.Range("E14").Value = "  +0.53%  "

$ws.Range("D15").Value = "59.371.19"
$ws.Range("E15").Value = "  +0.36%  "

$ws.Range("D16").Value = "21.04"
$ws.Range("E16").Value = "  +0.92%  "

$ws.Range("E17").Value = "  -0.81%  "

$ws.Range("D18").Value = "2.638.20"
$ws.Range("E18").Value = "  +0.20%  "

$ws.Range("D19").Value = "340.11"
$ws.Range("E19").Value = "  -2.39%  "

$ws.Range("E20").Value = "  -1.38%  "

$ws.Range("D21").Value = "10.31"
$ws.Range("E21").Value = "  +0.19%  "

$ws.Range("D22").Value = "6.28"
$ws.Range("E22").Value = "  +1.72%  "

$ws.Range("E23").Value = "  -0.29%  "

$ws.Range("E24").Value = "  +2.18%  "

$ws.Range("D25").Value = "0.167"
$ws.Range("E25").Value = "  +1.92%  "

$ws.Range("E26").Value = "  -0.79%  "

$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.48%  "

$ws.Range("D28").Value = "0.0₃0802"
$ws.Range("E28").Value = "  +0.03%  "

$ws.Range("E29").Value = "  +0.55%  "

$ws.Range("E30").Value = "  +4.67%  "

$ws.Range("E31").Value = "  +0.02%  "

$ws.Range("D32").Value = "1.59"
$ws.Range("E32").Value = "  +0.80%  "

$ws.Range("D33").Value = "18.72"
$ws.Range("E33").Value = "  -0.69%  "

$ws.Range("D34").Value = "149.49"
$ws.Range("E34").Value = "  +0.39%  "

$ws.Range("E35").Value = "  +3.09%  "

$ws.Range("D36").Value = "1.20"
$ws.Range("E36").Value = "  +2.22%  "

$ws.Range("E37").Value = "  -5.25%  "

$ws.Range("D38").Value = "0.883"
$ws.Range("E38").Value = "  +2.66%  "

$ws.Range("D39").Value = "36.93"
$ws.Range("E39").Value = "  +1.17%  "

$ws.Range("D40").Value = "1.49"
$ws.Range("E40").Value = "  +3.32%  "

$ws.Range("E41").Value = "  -1.70%  "

$ws.Range("E42").Value = "  +4.92%  "

$ws.Range("E43").Value = "  +0.15%  "

$ws.Range("D44").Value = "275.43"
$ws.Range("E44").Value = "  -0.62%  "

$ws.Range("D45").Value = "19.76"
$ws.Range("E45").Value = "  +1.09%  "

$ws.Range("D46").Value = "0.0974"
$ws.Range("E46").Value = "  -1.52%  "

$ws.Range("D47").Value = "0.0536"
$ws.Range("E47").Value = "  +1.73%  "

$ws.Range("D48").Value = "2.051.36"
$ws.Range("E48").Value = "  -1.96%  "

$ws.Range("E49").Value = "  +2.03%  "

$ws.Range("D50").Value = "4.77"
$ws.Range("E50").Value = "  +1.48%  "

$ws.Range("D51").Value = "0.0229"
$ws.Range("E51").Value = "  -0.77%  "
